$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: GX 460 Premium base MSRP bump
$ws.Range("D29").Value = 53100

# Row 30: GX 460 Luxury base MSRP bump
$ws.Range("D30").Value = 55890

# Row 31: LX 570 Two-Row base MSRP bump
$ws.Range("D31").Value = 64365

# Row 32: LX 570 Three-Row base MSRP bump + DPHF correction
$ws.Range("D32").Value = 86580
$ws.Range("E32").Value = 1025

# Row 33: LX 570 Inspiration Series SE base MSRP bump + DPHF correction
$ws.Range("D33").Value = 91580
$ws.Range("E33").Value = 1025

# Row 34: LX 570 Inspiration Series SE (second trim) -- was a placeholder text value,
# now filled in with the real base MSRP number, using the same number format as
# the surrounding MSRP column, plus the DPHF correction
$ws.Range("D34").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("D34").Value = 99310
$ws.Range("E34").Value = 1025
